$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 648:649, pushing the existing data (old rows 648-741)
# down to become rows 650-743.
$ws.Rows("648:649").Insert()

# Fill in the two new rows with the new "Región de O'Higgins" (guarda) records.
$ws.Range("A648").Value = 8
$ws.Range("B648").Value = "Terminal La Palmera de La Serena"
$ws.Range("C648").Value = "Coquimbo"
$ws.Range("D648").Value = 44776
$ws.Range("E648").Value = 4
$ws.Range("F648").Value = 100112045
$ws.Range("G648").Value = "Zapallo"
$ws.Range("H648").Value = "Camote"
$ws.Range("I648").Value = "1a (guarda)"
$ws.Range("J648").Value = 1800
$ws.Range("K648").Value = 1150
$ws.Range("L648").Value = 1200
$ws.Range("M648").Value = 1175
$ws.Range("N648").Value = "$/kilo (volumen en unidades)"
$ws.Range("O648").Value = "Región de O'Higgins"
$ws.Range("P648").Value = 1175
$ws.Range("Q648").Value = 1
$ws.Range("R648").Value = "Hortaliza"

$ws.Range("A649").Value = 8
$ws.Range("B649").Value = "Terminal La Palmera de La Serena"
$ws.Range("C649").Value = "Coquimbo"
$ws.Range("D649").Value = 44776
$ws.Range("E649").Value = 4
$ws.Range("F649").Value = 100112045
$ws.Range("G649").Value = "Zapallo"
$ws.Range("H649").Value = "Camote"
$ws.Range("I649").Value = "2a (guarda)"
$ws.Range("J649").Value = 960
$ws.Range("K649").Value = 1000
$ws.Range("L649").Value = 1050
$ws.Range("M649").Value = 1025
$ws.Range("N649").Value = "$/kilo (volumen en unidades)"
$ws.Range("O649").Value = "Región de O'Higgins"
$ws.Range("P649").Value = 1025
$ws.Range("Q649").Value = 1
$ws.Range("R649").Value = "Hortaliza"
